$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "57.931.05"
Set-TextValue "E2" "  +2.10%  "
Set-TextValue "D3" "3.058.60"
Set-TextValue "E3" "  +0.99%  "
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "516.52"
Set-TextValue "E5" "  +0.99%  "
Set-TextValue "D6" "141.73"
Set-TextValue "E6" "  +0.92%  "
Set-TextValue "E8" "  +1.06%  "
Set-TextValue "D9" "7.32"
Set-TextValue "E9" "  +2.62%  "
Set-TextValue "E10" "  -0.32%  "
Set-TextValue "D11" "0.375"
Set-TextValue "E11" "  +1.54%  "
Set-TextValue "D12" "3.579.98"
Set-TextValue "E12" "  +0.88%  "
Set-TextValue "E13" "  +3.02%  "
Set-TextValue "D14" "26.27"
Set-TextValue "E14" "  +3.66%  "
Set-TextValue "E15" "  +0.53%  "
Set-TextValue "D16" "57.946.38"
Set-TextValue "E16" "  +2.16%  "
Set-TextValue "D17" "3.051.38"
Set-TextValue "E17" "  +0.73%  "
Set-TextValue "D18" "6.11"
Set-TextValue "E18" "  +3.02%  "
Set-TextValue "D19" "12.81"
Set-TextValue "E19" "  -2.62%  "
Set-TextValue "D20" "8.06"
Set-TextValue "E20" "  +0.10%  "
Set-TextValue "D21" "331.73"
Set-TextValue "E21" "  -0.49%  "
Set-TextValue "D23" "0.501"
Set-TextValue "E23" "  -0.02%  "
Set-TextValue "D24" "65.46"
Set-TextValue "E24" "  +1.27%  "
Set-TextValue "E25" "  +2.79%  "
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  +0.21%  "
Set-TextValue "D27" "0.0₃0902"
Set-TextValue "E27" "  -2.93%  "
Set-TextValue "E28" "  +0.32%  "
Set-TextValue "E29" "  +6.44%  "
Set-TextValue "E30" "  +1.41%  "
Set-TextValue "E31" "  +3.14%  "
Set-TextValue "D32" "20.68"
Set-TextValue "E32" "  +1.28%  "
Set-TextValue "D33" "154.66"
Set-TextValue "E33" "  +1.30%  "
Set-TextValue "E34" "  +0.88%  "
Set-TextValue "D35" "5.99"
Set-TextValue "E35" "  +3.08%  "
Set-TextValue "D36" "26.92"
Set-TextValue "E36" "  -0.62%  "
Set-TextValue "D37" "1.27"
Set-TextValue "E37" "  +3.26%  "
Set-TextValue "D38" "0.0679"
Set-TextValue "E38" "  +2.52%  "
Set-TextValue "D39" "3.097.03"
Set-TextValue "E39" "  +0.99%  "
Set-TextValue "E40" "  +2.93%  "
Set-TextValue "D41" "36.58"
Set-TextValue "E41" "  +0.05%  "
Set-TextValue "E42" "  -0.02%  "
Set-TextValue "D43" "0.657"
Set-TextValue "E43" "  -0.61%  "
Set-TextValue "D44" "2.283.07"
Set-TextValue "E44" "  +3.15%  "
Set-TextValue "E45" "  +5.43%  "
Set-TextValue "E46" "  +1.62%  "
Set-TextValue "D47" "20.62"
Set-TextValue "E47" "  +4.11%  "
Set-TextValue "D48" "0.938"
Set-TextValue "E48" "  +0.54%  "
Set-TextValue "E49" "  +1.52%  "
Set-TextValue "D50" "0.733"
Set-TextValue "E50" "  +9.12%  "
Set-TextValue "D51" "0.0879"
Set-TextValue "E51" "  +2.73%  "
